$d = $word.ActiveDocument

# --- 1. Reset the "smarthosting" bookmark so it re-serializes with id=0 ---
# (Word assigns bookmark ids sequentially on save; deleting and re-adding
#  the bookmark over the same range causes it to be re-numbered starting
#  from 0, matching the target XML.)
$bm = $d.Bookmarks("smarthosting")
$bmRange = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("smarthosting", $bmRange)

# --- 2. Remove the CryptoBridge / HitBTC HYPERLINK fields entirely ---
# (begin/instrText/separate/result/end) using the Fields collection so the
# whole field structure disappears, not just the visible result text.
$cryptoBridgeField = $null
$hitBtcField = $null
foreach ($f in $d.Fields) {
    if ($f.Code.Text -like "*crypto-bridge.org*") {
        $cryptoBridgeField = $f
    } elseif ($f.Code.Text -like "*hitbtc.com*") {
        $hitBtcField = $f
    }
}
$hitBtcField.Delete()
$cryptoBridgeField.Delete()

# --- 3. Trim the leftover " such as" and ", " connector text ---
# (the source text uses non-breaking spaces (U+00A0) rather than regular
#  spaces around these words)
$nbsp = [char]0x00A0
$leftover = " such as" + $nbsp + "," + $nbsp
$null = $d.Content.Find.Execute($leftover, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
